$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 7).Value = 0.3227736666666667
$ws.Cells.Item(2, 8).Value = 0.968321
$ws.Cells.Item(2, 9).Value = 0.1416094457286952
$ws.Cells.Item(2, 10).Value = 0.1416094457286952
$ws.Cells.Item(2, 13).Value = 20.66830833333333
$ws.Cells.Item(2, 14).Value = 62.004925
$ws.Cells.Item(2, 15).Value = 0.6755285375771634
$ws.Cells.Item(2, 16).Value = 0.6755285375771634
$ws.Cells.Item(2, 17).Value = 6.671185664547222
$ws.Cells.Item(2, 18).Value = 60.040670980925
$ws.Cells.Item(2, 19).Value = 0.09566122178021812
$ws.Cells.Item(2, 20).Value = 0.09566122178021812

$ws.Cells.Item(3, 7).Value = 0.3227736666666667
$ws.Cells.Item(3, 8).Value = 0.968321
$ws.Cells.Item(3, 9).Value = 0.1416094457286952
$ws.Cells.Item(3, 10).Value = 0.1416094457286952
$ws.Cells.Item(3, 15).Value = 0.07047809033489469
$ws.Cells.Item(3, 16).Value = 0.07047809033489467
$ws.Cells.Item(3, 17).Value = 0.6960067558257779
$ws.Cells.Item(3, 18).Value = 6.264060802432
$ws.Cells.Item(3, 19).Value = 0.009980363308341345
$ws.Cells.Item(3, 20).Value = 0.009980363308341343

$ws.Cells.Item(4, 7).Value = 0.3227736666666667
$ws.Cells.Item(4, 8).Value = 0.968321
$ws.Cells.Item(4, 9).Value = 0.1416094457286952
$ws.Cells.Item(4, 10).Value = 0.1416094457286952
$ws.Cells.Item(4, 13).Value = 1.683564
$ws.Cells.Item(4, 14).Value = 5.050692
$ws.Cells.Item(4, 15).Value = 0.05502605769642779
$ws.Cells.Item(4, 16).Value = 0.05502605769642779
$ws.Cells.Item(4, 17).Value = 0.5434101253479999
$ws.Cells.Item(4, 18).Value = 4.890691128132
$ws.Cells.Item(4, 19).Value = 0.007792209531026341
$ws.Cells.Item(4, 20).Value = 0.007792209531026341

$ws.Cells.Item(5, 7).Value = 0.3227736666666667
$ws.Cells.Item(5, 8).Value = 0.968321
$ws.Cells.Item(5, 9).Value = 0.1416094457286952
$ws.Cells.Item(5, 10).Value = 0.1416094457286952
$ws.Cells.Item(5, 13).Value = 5.278649666666666
$ws.Cells.Item(5, 14).Value = 15.835949
$ws.Cells.Item(5, 15).Value = 0.1725288026574751
$ws.Cells.Item(5, 16).Value = 0.1725288026574751
$ws.Cells.Item(5, 17).Value = 1.703809107958778
$ws.Cells.Item(5, 18).Value = 15.334281971629
$ws.Cells.Item(5, 19).Value = 0.02443170811656048
$ws.Cells.Item(5, 20).Value = 0.02443170811656047

$ws.Cells.Item(6, 7).Value = 0.3227736666666667
$ws.Cells.Item(6, 8).Value = 0.968321
$ws.Cells.Item(6, 9).Value = 0.1416094457286952
$ws.Cells.Item(6, 10).Value = 0.1416094457286952
$ws.Cells.Item(6, 13).Value = 0.8089063333333334
$ws.Cells.Item(6, 14).Value = 2.426719
$ws.Cells.Item(6, 15).Value = 0.02643851173403914
$ws.Cells.Item(6, 16).Value = 0.02643851173403913
$ws.Cells.Item(6, 17).Value = 0.2610936631998889
$ws.Cells.Item(6, 18).Value = 2.349842968799
$ws.Cells.Item(6, 19).Value = 0.003743942992548885
$ws.Cells.Item(6, 20).Value = 0.003743942992548885

$ws.Cells.Item(7, 9).Value = 0.8226066833587575
$ws.Cells.Item(7, 10).Value = 0.8226066833587576
$ws.Cells.Item(7, 13).Value = 20.66830833333333
$ws.Cells.Item(7, 14).Value = 62.004925
$ws.Cells.Item(7, 15).Value = 0.6755285375771634
$ws.Cells.Item(7, 16).Value = 0.6755285375771634
$ws.Cells.Item(7, 17).Value = 38.75279565811945
$ws.Cells.Item(7, 18).Value = 348.775160923075
$ws.Cells.Item(7, 19).Value = 0.5556942898105421
$ws.Cells.Item(7, 20).Value = 0.5556942898105423

$ws.Cells.Item(8, 9).Value = 0.8226066833587575
$ws.Cells.Item(8, 10).Value = 0.8226066833587576
$ws.Cells.Item(8, 15).Value = 0.07047809033489469
$ws.Cells.Item(8, 16).Value = 0.07047809033489467
$ws.Cells.Item(8, 19).Value = 0.05797574813984663
$ws.Cells.Item(8, 20).Value = 0.05797574813984662

$ws.Cells.Item(9, 9).Value = 0.8226066833587575
$ws.Cells.Item(9, 10).Value = 0.8226066833587576
$ws.Cells.Item(9, 13).Value = 1.683564
$ws.Cells.Item(9, 14).Value = 5.050692
$ws.Cells.Item(9, 15).Value = 0.05502605769642779
$ws.Cells.Item(9, 16).Value = 0.05502605769642779
$ws.Cells.Item(9, 17).Value = 3.156659491292
$ws.Cells.Item(9, 18).Value = 28.409935421628
$ws.Cells.Item(9, 19).Value = 0.0452648028199661
$ws.Cells.Item(9, 20).Value = 0.04526480281996611

$ws.Cells.Item(10, 9).Value = 0.8226066833587575
$ws.Cells.Item(10, 10).Value = 0.8226066833587576
$ws.Cells.Item(10, 13).Value = 5.278649666666666
$ws.Cells.Item(10, 14).Value = 15.835949
$ws.Cells.Item(10, 15).Value = 0.1725288026574751
$ws.Cells.Item(10, 16).Value = 0.1725288026574751
$ws.Cells.Item(10, 17).Value = 9.897395983454556
$ws.Cells.Item(10, 18).Value = 89.07656385109101
$ws.Cells.Item(10, 19).Value = 0.1419233461379232
$ws.Cells.Item(10, 20).Value = 0.1419233461379232

$ws.Cells.Item(11, 9).Value = 0.8226066833587575
$ws.Cells.Item(11, 10).Value = 0.8226066833587576
$ws.Cells.Item(11, 13).Value = 0.8089063333333334
$ws.Cells.Item(11, 14).Value = 2.426719
$ws.Cells.Item(11, 15).Value = 0.02643851173403914
$ws.Cells.Item(11, 16).Value = 0.02643851173403913
$ws.Cells.Item(11, 17).Value = 1.516688319946778
$ws.Cells.Item(11, 18).Value = 13.650194879521
$ws.Cells.Item(11, 19).Value = 0.02174849645047953
$ws.Cells.Item(11, 20).Value = 0.02174849645047953

$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.081563
$ws.Cells.Item(12, 8).Value = 0.244689
$ws.Cells.Item(12, 9).Value = 0.03578387091254728
$ws.Cells.Item(12, 10).Value = 0.03578387091254728
$ws.Cells.Item(12, 13).Value = 20.66830833333333
$ws.Cells.Item(12, 14).Value = 62.004925
$ws.Cells.Item(12, 15).Value = 0.6755285375771634
$ws.Cells.Item(12, 16).Value = 0.6755285375771634
$ws.Cells.Item(12, 17).Value = 1.685769232591666
$ws.Cells.Item(12, 18).Value = 15.171923093325
$ws.Cells.Item(12, 19).Value = 0.02417302598640306
$ws.Cells.Item(12, 20).Value = 0.02417302598640306

$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.081563
$ws.Cells.Item(13, 8).Value = 0.244689
$ws.Cells.Item(13, 9).Value = 0.03578387091254728
$ws.Cells.Item(13, 10).Value = 0.03578387091254728
$ws.Cells.Item(13, 15).Value = 0.07047809033489469
$ws.Cells.Item(13, 16).Value = 0.07047809033489467
$ws.Cells.Item(13, 17).Value = 0.1758767981653333
$ws.Cells.Item(13, 18).Value = 1.582891183488
$ws.Cells.Item(13, 19).Value = 0.002521978886706717
$ws.Cells.Item(13, 20).Value = 0.002521978886706717

$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.081563
$ws.Cells.Item(14, 8).Value = 0.244689
$ws.Cells.Item(14, 9).Value = 0.03578387091254728
$ws.Cells.Item(14, 10).Value = 0.03578387091254728
$ws.Cells.Item(14, 13).Value = 1.683564
$ws.Cells.Item(14, 14).Value = 5.050692
$ws.Cells.Item(14, 15).Value = 0.05502605769642779
$ws.Cells.Item(14, 16).Value = 0.05502605769642779
$ws.Cells.Item(14, 17).Value = 0.137316530532
$ws.Cells.Item(14, 18).Value = 1.235848774788
$ws.Cells.Item(14, 19).Value = 0.001969045345435351
$ws.Cells.Item(14, 20).Value = 0.001969045345435351

$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.081563
$ws.Cells.Item(15, 8).Value = 0.244689
$ws.Cells.Item(15, 9).Value = 0.03578387091254728
$ws.Cells.Item(15, 10).Value = 0.03578387091254728
$ws.Cells.Item(15, 13).Value = 5.278649666666666
$ws.Cells.Item(15, 14).Value = 15.835949
$ws.Cells.Item(15, 15).Value = 0.1725288026574751
$ws.Cells.Item(15, 16).Value = 0.1725288026574751
$ws.Cells.Item(15, 17).Value = 0.4305425027623333
$ws.Cells.Item(15, 18).Value = 3.874882524860999
$ws.Cells.Item(15, 19).Value = 0.006173748402991432
$ws.Cells.Item(15, 20).Value = 0.006173748402991431

$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.081563
$ws.Cells.Item(16, 8).Value = 0.244689
$ws.Cells.Item(16, 9).Value = 0.03578387091254728
$ws.Cells.Item(16, 10).Value = 0.03578387091254728
$ws.Cells.Item(16, 13).Value = 0.8089063333333334
$ws.Cells.Item(16, 14).Value = 2.426719
$ws.Cells.Item(16, 15).Value = 0.02643851173403914
$ws.Cells.Item(16, 16).Value = 0.02643851173403913
$ws.Cells.Item(16, 17).Value = 0.06597682726566667
$ws.Cells.Item(16, 18).Value = 0.593791445391
$ws.Cells.Item(16, 19).Value = 0.0009460722910107228
$ws.Cells.Item(16, 20).Value = 0.0009460722910107227
